$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the paragraph that currently contains $searchText with the
# raw WordprocessingML supplied in $xml (a full <w:p>...</w:p> fragment).
# Uses Find to locate the anchor text so absolute character offsets never
# need to be hand-computed, then expands to the owning paragraph's Range.
# ---------------------------------------------------------------------------
function Replace-ParagraphContainingText($searchText, $xml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $searchText"
    }
    $para = $rng.Paragraphs(1)
    $target = $para.Range
    $target.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 1) End of document: after the "Data Analysis" paragraph, add a blank
#    paragraph, a "References" Heading2 paragraph, and another blank
#    paragraph (before the trailing sectPr).
# ---------------------------------------------------------------------------
$rngEnd = $d.Content
$foundEnd = $rngEnd.Find.Execute("Data Analysis", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundEnd) {
    throw "Could not find text: Data Analysis"
}
$dataAnalysisPara = $rngEnd.Paragraphs(1)
$endPos = $dataAnalysisPara.Range.End
$endInsertPoint = $d.Range($endPos, $endPos)
$endXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2"/><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>References</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr></w:pPr></w:p>
'@
$endInsertPoint.InsertXML($endXml)

# ---------------------------------------------------------------------------
# 2) "Goes from 3-d point cloud to joint/finger estimation" paragraph: split
#    the final word off into its own (grammar-flagged) run, AND delete the
#    eight list paragraphs that used to follow it (Hierarchical RL .. RL
#    Games Conference), all the way up to (not including) the
#    "** I will require computing power" paragraph.
# ---------------------------------------------------------------------------
$rng45 = $d.Content
$found45 = $rng45.Find.Execute("Goes from 3-d point cloud to joint/finger estimation", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found45) {
    throw "Could not find text: Goes from 3-d point cloud to joint/finger estimation"
}
$startPara = $rng45.Paragraphs(1)

$rngEndAnchor = $d.Content
$foundEndAnchor = $rngEndAnchor.Find.Execute("RL Games Conference", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundEndAnchor) {
    throw "Could not find text: RL Games Conference"
}
$endPara = $rngEndAnchor.Paragraphs(1)

$span = $d.Range($startPara.Range.Start, $endPara.Range.End)
$p45Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="12"/></w:numPr><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr><w:t xml:space="preserve">Goes from 3-d point cloud to joint/finger </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>estimation</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>
'@
$span.InsertXML($p45Xml)

# ---------------------------------------------------------------------------
# 3) "Proposes feature extraction using CNN, along with an RL module for
#    path optimisation" -> split off "optimisation" with gramStart/gramEnd.
# ---------------------------------------------------------------------------
$p44Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="12"/></w:numPr><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr><w:t xml:space="preserve">Proposes feature extraction using CNN, along with an RL module for path </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>optimisation</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>
'@
Replace-ParagraphContainingText "Proposes feature extraction using CNN, along with an RL module for path optimisation" $p44Xml

# ---------------------------------------------------------------------------
# 4) "Very informative on concepts within RL and gives examples with code"
#    -> add a lastRenderedPageBreak and split off "code" with gramStart/End.
# ---------------------------------------------------------------------------
$p42Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="12"/></w:numPr><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Very informative on concepts within RL and gives examples with </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>code</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>
'@
Replace-ParagraphContainingText "Very informative on concepts within RL and gives examples with code" $p42Xml

# ---------------------------------------------------------------------------
# 5) "Textbook on Reinforcememnt learning with Python" -> split off "Python"
#    with gramStart/gramEnd.
# ---------------------------------------------------------------------------
$p41Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="12"/></w:numPr><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr><w:t xml:space="preserve">Textbook on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>Reinforcememnt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr><w:t xml:space="preserve"> learning with </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr><w:t>Python</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>
'@
Replace-ParagraphContainingText "Textbook on Reinforcememnt learning with Python" $p41Xml

# ---------------------------------------------------------------------------
# 6) Fill the empty paragraph directly under the "Reading" heading with the
#    new write-up paragraph about AAAI papers / imitation learning / RL.
# ---------------------------------------------------------------------------
$rngReading = $d.Content
$foundReading = $rngReading.Find.Execute("Reading", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundReading) {
    throw "Could not find text: Reading"
}
$headingPara = $rngReading.Paragraphs(1)
$targetPara = $headingPara.Next()
$para37Xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="en-GB"/></w:rPr><w:t xml:space="preserve">Many papers for the AAAI Conference propose RL strategies for games. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Gema</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Parreño Piqueras</w:t></w:r><w:r><w:t xml:space="preserve"> [] proposes an Imitation Learning approach for game designers wanting to develop a companion style AI</w:t></w:r><w:r><w:t xml:space="preserve">, which requires a training phase to learn a policy that mimics behaviour of an expert. Self-driving cars also need systems to drive effectively, J. Duan et al proposed a hierarchical RL approach, where a manoeuvre policy was used, giving a manoeuvre based on the environment, which was then fed into a sub-policy detailing what brake/accelerator </w:t></w:r><w:r><w:t>and steering inputs are needed. This had better results when compared to normal RL, such as shorter training time and higher average reward. I considered taking this approach, specifying a set of inputs resulting in; sharp turn left/right, slight turn left/right and straight on, but decided against it as I would prefer to see the AI learn these things by itself.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p>
'@
$targetPara.Range.InsertXML($para37Xml)

Write-Output "Done"
